$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.370.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +9.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.675.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9995"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.40%  "

$ws.Range("E7").Value = "  +2.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3417"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +14.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.156"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07205"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.02%  "

$ws.Range("E12").Value = "  -0.46%  "

$ws.Range("E13").Value = "  +5.59%  "

$ws.Range("E14").Value = "  +4.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.718"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.674.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9986"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06645"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "80.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.086"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.312.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.435"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.26%  "

$ws.Range("E26").Value = "  +7.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.861.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.263"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.045"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9672"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08438"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.675"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06384"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.304"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02311"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.658"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.228"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2081"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6063"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9992"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("E45").Value = "  -0.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5854"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.014"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07141"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.09%  "
